$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. G6 becomes a new "Next" marker cell.
#    H4 currently carries the "Next" marker's look (red font on the blue fill),
#    so clone that formatting onto G6 before H4's own value/format change below.
$ws.Range("H4").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "Next"

# 2. H4's "Next" placeholder resolves to an actual date; it should now look
#    like the other resolved date cells in that column/fill group (e.g. G4).
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 45868

# 3. O7 gets filled in with a resolved date; its existing formatting is kept.
$ws.Range("O7").Value = 45868

$excel.CutCopyMode = 0

# 4. Update the active selection to match.
$ws.Range("R17").Select()
